$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)

# The paragraph currently has 3 separate runs ("An", " ", "image").
# Re-assigning TextRange.Text with the same rendered string is a no-op in
# this engine's diffing, so first set a distinct placeholder to force a
# real text rewrite, then set the final text -- this collapses the
# paragraph down to a single run containing "An image".
$sh.TextFrame.TextRange.Text = "__placeholder__"
$sh.TextFrame.TextRange.Text = "An image"
